$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4378.9287
$ws.Cells.Item(62, 9).Value = 2609.9
$ws.Cells.Item(62, 11).Value = 2609.9
$ws.Cells.Item(62, 13).Value = -1985.9
$ws.Cells.Item(65, 8).Value = 4378.9287
$ws.Cells.Item(65, 9).Value = 2609.9
$ws.Cells.Item(65, 11).Value = 13049.5
$ws.Cells.Item(65, 13).Value = -9929.5
$ws.Cells.Item(70, 8).Value = 13342604
$ws.Cells.Item(70, 10).Value = 20010500
$ws.Cells.Item(70, 12).Value = 60031500
$ws.Cells.Item(70, 14).Value = -60032040
$ws.Cells.Item(73, 8).Value = 13342604
$ws.Cells.Item(73, 10).Value = 20010500
$ws.Cells.Item(73, 12).Value = 60031500
$ws.Cells.Item(73, 14).Value = -60033372
$ws.Cells.Item(100, 8).Value = 4572.7144
$ws.Cells.Item(100, 9).Value = 2480.6
$ws.Cells.Item(100, 11).Value = 2480.6
$ws.Cells.Item(100, 13).Value = -1939.6
$ws.Cells.Item(112, 8).Value = 1324.45
$ws.Cells.Item(112, 10).Value = 1338.2778
$ws.Cells.Item(112, 12).Value = 4014.8334
$ws.Cells.Item(112, 14).Value = -6230.8334
$ws.Cells.Item(132, 8).Value = 1732.119
$ws.Cells.Item(132, 9).Value = 1648.0769
$ws.Cells.Item(132, 10).Value = 2824.6667
$ws.Cells.Item(132, 11).Value = 4944.2307
$ws.Cells.Item(132, 12).Value = 8474.000100000001
$ws.Cells.Item(132, 13).Value = -2414.2307
$ws.Cells.Item(132, 14).Value = -13534.0001
$ws.Cells.Item(135, 8).Value = 3558
$ws.Cells.Item(135, 9).Value = 2933.3333
$ws.Cells.Item(135, 11).Value = 26399.9997
$ws.Cells.Item(135, 13).Value = -23864.9997
$ws.Cells.Item(138, 8).Value = 2444.5366
$ws.Cells.Item(138, 9).Value = 1108.9375
$ws.Cells.Item(138, 10).Value = 3299.32
$ws.Cells.Item(138, 11).Value = 3326.8125
$ws.Cells.Item(138, 12).Value = 9897.960000000001
$ws.Cells.Item(138, 13).Value = 1813.1875
$ws.Cells.Item(138, 14).Value = -20177.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 10550.75
$ws.Cells.Item(132, 9).Value = 10062.75
$ws.Cells.Item(132, 10).Value = 12502.75
$ws.Cells.Item(132, 11).Value = 30188.25
$ws.Cells.Item(132, 12).Value = 37508.25
$ws.Cells.Item(132, 13).Value = -27658.25
$ws.Cells.Item(132, 14).Value = -42568.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 7043
$ws.Cells.Item(86, 9).Value = 5547.5
$ws.Cells.Item(86, 11).Value = 5547.5
$ws.Cells.Item(86, 13).Value = -4424.5
$ws.Cells.Item(89, 8).Value = 7043
$ws.Cells.Item(89, 9).Value = 5547.5
$ws.Cells.Item(89, 11).Value = 27737.5
$ws.Cells.Item(89, 13).Value = -22121.5
$ws.Cells.Item(94, 8).Value = 1799.4286
$ws.Cells.Item(94, 9).Value = 1287.0769
$ws.Cells.Item(94, 10).Value = 2632
$ws.Cells.Item(94, 11).Value = 1287.0769
$ws.Cells.Item(94, 12).Value = 2632
$ws.Cells.Item(94, 13).Value = -836.0769
$ws.Cells.Item(94, 14).Value = -3534
$ws.Cells.Item(105, 8).Value = 17238.889
$ws.Cells.Item(105, 9).Value = 21690.2
$ws.Cells.Item(105, 11).Value = 21690.2
$ws.Cells.Item(105, 13).Value = -19943.2
$ws.Cells.Item(134, 8).Value = 1499.7778
$ws.Cells.Item(134, 9).Value = 1499.7778
$ws.Cells.Item(134, 11).Value = 4499.3334
$ws.Cells.Item(134, 13).Value = -1964.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 14096.8
$ws.Cells.Item(62, 9).Value = 4159.3335
$ws.Cells.Item(62, 10).Value = 29003
$ws.Cells.Item(62, 11).Value = 4159.3335
$ws.Cells.Item(62, 12).Value = 29003
$ws.Cells.Item(62, 13).Value = -3535.3335
$ws.Cells.Item(62, 14).Value = -30251
$ws.Cells.Item(65, 8).Value = 14096.8
$ws.Cells.Item(65, 9).Value = 4159.3335
$ws.Cells.Item(65, 10).Value = 29003
$ws.Cells.Item(65, 11).Value = 20796.6675
$ws.Cells.Item(65, 12).Value = 145015
$ws.Cells.Item(65, 13).Value = -17676.6675
$ws.Cells.Item(65, 14).Value = -151255
$ws.Cells.Item(94, 8).Value = 2593.7
$ws.Cells.Item(94, 9).Value = 736
$ws.Cells.Item(94, 10).Value = 3058.125
$ws.Cells.Item(94, 11).Value = 736
$ws.Cells.Item(94, 12).Value = 3058.125
$ws.Cells.Item(94, 13).Value = -285
$ws.Cells.Item(94, 14).Value = -3960.125
$ws.Cells.Item(131, 8).Value = 70000
$ws.Cells.Item(131, 10).Value = 70000
$ws.Cells.Item(131, 12).Value = 70000
$ws.Cells.Item(131, 14).Value = -80080
$ws.Cells.Item(132, 8).Value = 3213.923
$ws.Cells.Item(132, 9).Value = 2269.9443
$ws.Cells.Item(132, 10).Value = 5337.875
$ws.Cells.Item(132, 11).Value = 6809.8329
$ws.Cells.Item(132, 12).Value = 16013.625
$ws.Cells.Item(132, 13).Value = -4279.8329
$ws.Cells.Item(132, 14).Value = -21073.625
$ws.Cells.Item(134, 8).Value = 4228.8823
$ws.Cells.Item(134, 9).Value = 2634.1428
$ws.Cells.Item(134, 10).Value = 11671
$ws.Cells.Item(134, 11).Value = 7902.428400000001
$ws.Cells.Item(134, 12).Value = 35013
$ws.Cells.Item(134, 13).Value = -5367.428400000001
$ws.Cells.Item(134, 14).Value = -40083

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(21, 8).Value = 180.27272
$ws.Cells.Item(21, 9).Value = 40.42857
$ws.Cells.Item(21, 10).Value = 425
$ws.Cells.Item(21, 11).Value = 121.28571
$ws.Cells.Item(21, 12).Value = 1275
$ws.Cells.Item(21, 13).Value = 51.71429000000001
$ws.Cells.Item(21, 14).Value = -1621
$ws.Cells.Item(92, 8).Value = 2807.9443
$ws.Cells.Item(92, 10).Value = 3203.3076
$ws.Cells.Item(92, 12).Value = 9609.9228
$ws.Cells.Item(92, 14).Value = -12105.9228
$ws.Cells.Item(122, 8).Value = 1300.5555
$ws.Cells.Item(122, 10).Value = 2172.1667
$ws.Cells.Item(122, 12).Value = 19549.5003
$ws.Cells.Item(122, 14).Value = -24449.5003

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5228.077
$ws.Cells.Item(7, 9).Value = 4326.575
$ws.Cells.Item(7, 10).Value = 8233.083000000001
$ws.Cells.Item(7, 11).Value = 4326.575
$ws.Cells.Item(7, 12).Value = 8233.083000000001
$ws.Cells.Item(7, 13).Value = -4214.575
$ws.Cells.Item(7, 14).Value = -8457.083000000001
$ws.Cells.Item(22, 8).Value = 2860.9565
$ws.Cells.Item(22, 9).Value = 1100.25
$ws.Cells.Item(22, 11).Value = 1100.25
$ws.Cells.Item(22, 13).Value = -805.25
$ws.Cells.Item(27, 8).Value = 2860.9565
$ws.Cells.Item(27, 9).Value = 1100.25
$ws.Cells.Item(27, 11).Value = 1100.25
$ws.Cells.Item(27, 13).Value = -993.25
$ws.Cells.Item(46, 8).Value = 3991.4119
$ws.Cells.Item(46, 9).Value = 3625
$ws.Cells.Item(46, 10).Value = 4104.154
$ws.Cells.Item(46, 11).Value = 3625
$ws.Cells.Item(46, 12).Value = 4104.154
$ws.Cells.Item(46, 13).Value = -3437
$ws.Cells.Item(46, 14).Value = -4480.154
$ws.Cells.Item(61, 8).Value = 4515.7896
$ws.Cells.Item(61, 9).Value = 3113.5
$ws.Cells.Item(61, 11).Value = 3113.5
$ws.Cells.Item(61, 13).Value = -2911.5
$ws.Cells.Item(68, 8).Value = 4535.091
$ws.Cells.Item(68, 9).Value = 3342.1428
$ws.Cells.Item(68, 10).Value = 6622.75
$ws.Cells.Item(68, 11).Value = 3342.1428
$ws.Cells.Item(68, 12).Value = 6622.75
$ws.Cells.Item(68, 13).Value = -2593.1428
$ws.Cells.Item(68, 14).Value = -8120.75
$ws.Cells.Item(71, 8).Value = 4535.091
$ws.Cells.Item(71, 9).Value = 3342.1428
$ws.Cells.Item(71, 10).Value = 6622.75
$ws.Cells.Item(71, 11).Value = 16710.714
$ws.Cells.Item(71, 12).Value = 33113.75
$ws.Cells.Item(71, 13).Value = -12966.714
$ws.Cells.Item(71, 14).Value = -40601.75
$ws.Cells.Item(93, 8).Value = 1443.6522
$ws.Cells.Item(93, 9).Value = 1289.7778
$ws.Cells.Item(93, 11).Value = 1289.7778
$ws.Cells.Item(93, 13).Value = -41.77780000000007
$ws.Cells.Item(100, 8).Value = 9189.25
$ws.Cells.Item(100, 9).Value = 2751.5
$ws.Cells.Item(100, 10).Value = 11335.167
$ws.Cells.Item(100, 11).Value = 2751.5
$ws.Cells.Item(100, 12).Value = 11335.167
$ws.Cells.Item(100, 13).Value = -2210.5
$ws.Cells.Item(100, 14).Value = -12417.167
$ws.Cells.Item(113, 8).Value = 4515.7896
$ws.Cells.Item(113, 9).Value = 3113.5
$ws.Cells.Item(113, 11).Value = 3113.5
$ws.Cells.Item(113, 13).Value = -943.5
$ws.Cells.Item(126, 8).Value = 5228.077
$ws.Cells.Item(126, 9).Value = 4326.575
$ws.Cells.Item(126, 10).Value = 8233.083000000001
$ws.Cells.Item(126, 11).Value = 12979.725
$ws.Cells.Item(126, 12).Value = 24699.249
$ws.Cells.Item(126, 13).Value = -10509.725
$ws.Cells.Item(126, 14).Value = -29639.249
$ws.Cells.Item(136, 8).Value = 15021.777
$ws.Cells.Item(136, 10).Value = 15649.5
$ws.Cells.Item(136, 12).Value = 46948.5
$ws.Cells.Item(136, 14).Value = -52048.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 9542.223
$ws.Cells.Item(62, 10).Value = 11141.857
$ws.Cells.Item(62, 12).Value = 11141.857
$ws.Cells.Item(62, 14).Value = -12389.857
$ws.Cells.Item(65, 8).Value = 9542.223
$ws.Cells.Item(65, 10).Value = 11141.857
$ws.Cells.Item(65, 12).Value = 55709.285
$ws.Cells.Item(65, 14).Value = -61949.285
$ws.Cells.Item(113, 8).Value = 525.4074000000001
$ws.Cells.Item(113, 9).Value = 384.6154
$ws.Cells.Item(113, 10).Value = 656.1429000000001
$ws.Cells.Item(113, 11).Value = 1153.8462
$ws.Cells.Item(113, 12).Value = 1968.4287
$ws.Cells.Item(113, 13).Value = 1016.1538
$ws.Cells.Item(113, 14).Value = -6308.4287
$ws.Cells.Item(126, 8).Value = 1861.5151
$ws.Cells.Item(126, 9).Value = 1645.963
$ws.Cells.Item(126, 11).Value = 4937.889
$ws.Cells.Item(126, 13).Value = -2467.889
$ws.Cells.Item(132, 8).Value = 13000.8
$ws.Cells.Item(132, 9).Value = 6999.5
$ws.Cells.Item(132, 11).Value = 20998.5
$ws.Cells.Item(132, 13).Value = -18468.5
$ws.Cells.Item(136, 8).Value = 2443.4092
$ws.Cells.Item(136, 9).Value = 2071.5715
$ws.Cells.Item(136, 11).Value = 6214.7145
$ws.Cells.Item(136, 13).Value = -3664.7145
